$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = 'Datos actualizados a 21 de Mayo de 2020 a las 06:05'

# Row 4
$ws.Range("B4").Value = 1592723
$ws.Range("E4").Value = 1127711
$ws.Range("H4").Value = 94936

# Row 14
$ws.Range("B14").Value = 112359
$ws.Range("C14").Value = 331
$ws.Range("E14").Value = 63502
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 3435

# Row 22
$ws.Range("B22").Value = 48091
$ws.Range("C22").Value = 2193
$ws.Range("D22").Value = 14155
$ws.Range("E22").Value = 32919
$ws.Range("G22").Value = 32
$ws.Range("H22").Value = 1017

# Row 42
$ws.Range("B42").Value = 16385
$ws.Range("D42").Value = 12286
$ws.Range("E42").Value = 3328
$ws.Range("H42").Value = 771

# Row 57
$ws.Range("A57").Value = 'Kazajistan'
$ws.Range("B57").Value = 7234
$ws.Range("C57").Value = 265
$ws.Range("D57").Value = 3734
$ws.Range("E57").Value = 3465
$ws.Range("H57").Value = 35

# Row 58
$ws.Range("A58").Value = 'Marruecos'
$ws.Range("B58").Value = 7133
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 4098
$ws.Range("E58").Value = 2841
$ws.Range("H58").Value = 194

# Row 59
$ws.Range("A59").Value = 'Australia'
$ws.Range("B59").Value = 7081
$ws.Range("C59").Value = 2
$ws.Range("D59").Value = 6470
$ws.Range("E59").Value = 511
$ws.Range("H59").Value = 100

# Row 60
$ws.Range("A60").Value = 'Malasia'
$ws.Range("B60").Value = 7009
$ws.Range("D60").Value = 5706
$ws.Range("E60").Value = 1189
$ws.Range("H60").Value = 114

# Row 73
$ws.Range("A73").Value = 'Sudan'
$ws.Range("B73").Value = 3138
$ws.Range("C73").Value = 410
$ws.Range("D73").Value = 309
$ws.Range("E73").Value = 2708
$ws.Range("G73").Value = 10
$ws.Range("H73").Value = 121

# Row 74
$ws.Range("B74").Value = 3100
$ws.Range("C74").Value = 145
$ws.Range("D74").Value = 355
$ws.Range("E74").Value = 2594
$ws.Range("G74").Value = 4
$ws.Range("H74").Value = 151

# Row 75
$ws.Range("A75").Value = 'Tailandia'
$ws.Range("B75").Value = 3034
$ws.Range("D75").Value = 2888
$ws.Range("E75").Value = 90
$ws.Range("H75").Value = 56

# Row 76
$ws.Range("A76").Value = 'Uzbekistan'
$ws.Range("B76").Value = 2939
$ws.Range("D76").Value = 2372
$ws.Range("E76").Value = 554
$ws.Range("H76").Value = 13

# Row 77
$ws.Range("A77").Value = 'Guinea'
$ws.Range("B77").Value = 2863
$ws.Range("D77").Value = 1525
$ws.Range("E77").Value = 1320
$ws.Range("H77").Value = 18

# Row 78
$ws.Range("A78").Value = 'Grecia'
$ws.Range("B78").Value = 2850
$ws.Range("D78").Value = 1374
$ws.Range("E78").Value = 1310
$ws.Range("H78").Value = 166

# Row 94
$ws.Range("E94").Value = 1002
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 32

# Row 143
$ws.Range("D143").Value = 264
$ws.Range("E143").Value = 60

# Row 197
$ws.Range("A197").Value = 'Nueva Caledonia'
$ws.Range("D197").Value = 18
$ws.Range("H197").Value = 0

# Row 199
$ws.Range("A199").Value = 'Belice'
$ws.Range("D199").Value = 16
$ws.Range("H199").Value = 2

# Row 209
$ws.Range("A209").Value = 'Seychelles'

# Row 211
$ws.Range("A211").Value = 'Groenlandia'

# Row 214
$ws.Range("A214").Value = 'Sahara Occidental'

# Row 215
$ws.Range("A215").Value = 'Bonaire, San Eustaquio y Saba'
